$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.662.04"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "3.777.47"
$ws.Range("E3").Value = "  -1.85%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.70"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("D7").Value = "3.774.97"
$ws.Range("E7").Value = "  -1.91%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -0.94%  "

$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.452"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000276"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.31%  "

$ws.Range("D15").Value = "4.412.77"
$ws.Range("E15").Value = "  -1.90%  "

$ws.Range("D16").Value = "3.780.38"
$ws.Range("E16").Value = "  -1.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("D18").Value = "67.630.98"
$ws.Range("E18").Value = "  -1.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.65%  "

$ws.Range("E20").Value = "  +0.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.716"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.19%  "

$ws.Range("E24").Value = "  -7.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("E30").Value = "  -1.72%  "

$ws.Range("D31").Value = "3.930.37"
$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.67%  "

$ws.Range("D36").Value = "3.743.15"
$ws.Range("E36").Value = "  -1.87%  "

$ws.Range("E37").Value = "  -1.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "

$ws.Range("E39").Value = "  -1.72%  "

$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.310"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.41%  "

$ws.Range("E44").Value = "  +0.00%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "

$ws.Range("E46").Value = "  -2.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "395.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000268"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.98%  "

Write-Host "Applied cryptos update"